$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 8, leaving only the first three rows
$ws.Range("A4:A8").EntireRow.Delete() | Out-Null

# Update the remaining cell values to the new content
$ws.Range("A1").Value = "URL"
$ws.Range("A2").Value = "https://www.rybelsus.com/"
$ws.Range("A3").Value = "https://www.rybelsus.com/savings-and-support.html"

# Restore the selected cell as shown in the updated sheet view
$ws.Range("A2").Select() | Out-Null
